# Auto-generated edit script: update crypto price/volume table
# to match the scraped values committed on Tue Jan 30 14:33:55 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.366.20"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "2.312.51"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'309.02"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "'105.33"
$ws.Range("E6").Value = "  +9.29%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +6.61%  "
$ws.Range("D10").Value = "'35.89"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").Value = "'52.66"
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "'6.98"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "2.668.64"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "'15.20"
$ws.Range("E16").Value = "  +5.94%  "
$ws.Range("D17").Value = "2.310.64"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "'0.802"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "43.309.73"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").Value = "'11.94"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("E22").Value = "  +4.98%  "
$ws.Range("D23").Value = "'67.93"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("D24").Value = "'240.56"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("D26").Value = "'2.61"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "'24.77"
$ws.Range("E28").Value = "  +6.88%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'36.22"
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.59"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.11"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "'162.20"
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'18.32"
$ws.Range("E35").Value = "  +5.46%  "
$ws.Range("E36").Value = "  +6.48%  "
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.57"
$ws.Range("E38").Value = "  +13.52%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'3.01"
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("E40").Value = "  +3.62%  "
$ws.Range("E41").Value = "  +4.37%  "
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").Value = "'2.52"
$ws.Range("E43").Value = "  +16.17%  "
$ws.Range("E44").Value = "  +3.41%  "
$ws.Range("D45").Value = "1.968.30"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "'18.75"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").Value = "'3.08"
$ws.Range("E47").Value = "  +6.93%  "
$ws.Range("D48").Value = "'10.24"
$ws.Range("E48").Value = "  +6.48%  "
$ws.Range("D49").Value = "'57.98"
$ws.Range("E49").Value = "  +8.07%  "
$ws.Range("D50").Value = "'2.95"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").Value = "'1.59"
$ws.Range("E51").Value = "  +9.14%  "
